# Automatische test-sync: 2025-06-18 12:30:10
# Appends the new incoming-mail log entry to the "Logs" sheet and refreshes
# the category counts on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append row 10 -------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A10").Value = "Vragen over samenwerking"
$logs.Range("B10").Value = "mailmind.test@zohomail.eu"
$logs.Range("C10").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D10").Value = "Overig"
$logs.Range("F10").Value = "2025-06-18 12:00:11"
$logs.Range("G10").Value = "Nee"

# --- "Dashboard" sheet: refresh category counts -----------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Overig"
$dash.Range("B3").Value = 2

$dash.Range("A4").Value = "Afmelding"
$dash.Range("B4").Value = 2

# --- "Logs" sheet: extend conditional formatting to include the new row ----------
$catFcs = $logs.Range("D2:D9").FormatConditions
for ($i = 1; $i -le $catFcs.Count; $i++) {
    $catFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D10"))
}

$answeredFcs = $logs.Range("G2:G9").FormatConditions
for ($i = 1; $i -le $answeredFcs.Count; $i++) {
    $answeredFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G10"))
}
